# Refresh cryptos price/volume snapshot (values + 2 swapped rows: WETH/WBTC)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.773.16"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "'3.415.03"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'596.12"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'141.69"
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("D7").Value = "'3.414.71"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "'8.00"
$ws.Range("E10").Value = "  +5.83%  "
$ws.Range("E11").Value = "  -5.70%  "
$ws.Range("D12").Value = "'0.406"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").Value = "'3.993.53"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "'0.0000199"
$ws.Range("E14").Value = "  -5.96%  "
$ws.Range("D15").Value = "'29.61"
$ws.Range("E15").Value = "  -5.04%  "
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'65.794.40"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.410.88"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").Value = "'10.30"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("E20").Value = "  -4.75%  "
$ws.Range("D21").Value = "'14.53"
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("D22").Value = "'414.51"
$ws.Range("E22").Value = "  -4.36%  "
$ws.Range("D23").Value = "'0.575"
$ws.Range("E23").Value = "  -5.15%  "
$ws.Range("D24").Value = "'77.17"
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -8.30%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("E27").Value = "  -4.94%  "
$ws.Range("D28").Value = "'7.92"
$ws.Range("E28").Value = "  -5.44%  "
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("E32").Value = "  -8.17%  "
$ws.Range("D33").Value = "'24.59"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").Value = "'3.412.40"
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  -6.22%  "
$ws.Range("D37").Value = "'5.48"
$ws.Range("E37").Value = "  -6.99%  "
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'168.80"
$ws.Range("E40").Value = "  -3.65%  "
$ws.Range("D41").Value = "'0.0851"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("D43").Value = "'5.03"
$ws.Range("E43").Value = "  -6.68%  "
$ws.Range("D44").Value = "'1.89"
$ws.Range("E44").Value = "  -10.04%  "
$ws.Range("D45").Value = "'45.37"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("D46").Value = "'26.26"
$ws.Range("E46").Value = "  -8.67%  "
$ws.Range("D47").Value = "'1.20"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("D48").Value = "'7.05"
$ws.Range("E48").Value = "  -5.00%  "
$ws.Range("D49").Value = "'2.28"
$ws.Range("E49").Value = "  -6.07%  "
$ws.Range("D50").Value = "'0.918"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").Value = "'0.231"
$ws.Range("E51").Value = "  -5.56%  "
